$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 109, shifting existing rows (and their
# weekly price records) down by one.
$ws.Rows.Item(109).Insert()

# Populate the newly inserted row 109 with the new weekly record.
$ws.Range("A109").Value = 11
$ws.Range("B109").Value = "Vega Monumental Concepción"
$ws.Range("C109").Value = "Bíobío"
$ws.Range("D109").Value = 44601
$ws.Range("E109").Value = 8
$ws.Range("F109").Value = 100112023
$ws.Range("G109").Value = "Brócoli"
$ws.Range("H109").Value = "Sin especificar"
$ws.Range("I109").Value = "Primera"
$ws.Range("J109").Value = 1500
$ws.Range("K109").Value = 750
$ws.Range("L109").Value = 800
$ws.Range("M109").Value = 773
$ws.Range("N109").Value = "$/unidad"
$ws.Range("O109").Value = "Región del Maule"
$ws.Range("P109").Value = 773
$ws.Range("Q109").Value = 1
$ws.Range("R109").Value = "Hortaliza"
